# settings.xlsx update
# - do_normalization (B7): 1 -> 0
# - do_database_injection (B15): 1 -> 0
# - do_subsetting (B17): 1 -> 0
# - data_subsets (B24): reorder the comma-separated list
# - selection/active cell moves to B24 (was C24), view scrolled down a bit

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Flip the three boolean toggles off.
$ws.Range("B7").Value = 0
$ws.Range("B15").Value = 0
$ws.Range("B17").Value = 0

# Re-order the data_subsets value.
$ws.Range("B24").Value = "CD4_T, CD8_T, TCRgd_T, B, Monos_and_DCs, NK"

# Match the new selection/scroll position recorded in the saved view.
$aw = $excel.ActiveWindow
$aw.ScrollRow = 19
$aw.ScrollColumn = 1
$ws.Range("B24").Select()

"settings updated"
